$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the sign-up email/password test data:
# A2 previously held "test@yopmail.com" -> now "yes@shuramail.com"
# B2 keeps "Test$1234!"
$ws.Range("A2").Value = "yes@shuramail.com"
$ws.Range("B2").Value = 'Test$1234!'

$wb.Save()
